$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A:T) reflecting the refreshed TPM-based NATMI output.
# String columns (A,B,C,D) use cluster/gene names; numeric columns are plain doubles.
$data = @(
    @("ECs",  "Timp3", "Agtr2", "ECs",  3, 1, 125.60613,          376.81839,          0.5776575759508482, 0.5776575759508483, 1, 0.3333333333333333, 0.045282,           0.135846,  0.001264082841858775, 0.001264082841858776, 5.687696778659999,  51.18927100794,     0.0007302070302291995, 0.0007302070302291999),
    @("ECs",  "Timp3", "Agtr2", "FAPs", 3, 1, 125.60613,          376.81839,          0.5776575759508482, 0.5776575759508483, 3, 1,                   34.682839,          104.048517, 0.9681988800594137,  0.9681988800594138,  4356.37718420307,   39207.39465782764,  0.559287418093447,    0.5592874180934471),
    @("ECs",  "Timp3", "Agtr2", "MuSCs",3, 1, 125.60613,          376.81839,          0.5776575759508482, 0.5776575759508483, 2, 0.6666666666666666, 1.093898333333333, 3.281695,  0.03053703709872749, 0.03053703709872749, 137.40033626345,    1236.60302637105,   0.01763995082717204,  0.01763995082717205),
    @("FAPs", "Timp3", "Agtr2", "ECs",  3, 1, 36.46294533333333, 109.388836,         0.1676916294872044, 0.1676916294872044, 1, 0.3333333333333333, 0.045282,           0.135846,  0.001264082841858775, 0.001264082841858776, 1.651115090584,     14.860035815256,    0.0002119761115581141,0.0002119761115581142),
    @("FAPs", "Timp3", "Agtr2", "FAPs", 3, 1, 36.46294533333333, 109.388836,         0.1676916294872044, 0.1676916294872044, 3, 1,                   34.682839,          104.048517, 0.9681988800594137,  0.9681988800594138,  1264.638462461801,  11381.74616215621,  0.1623588478648494,   0.1623588478648495),
    @("FAPs", "Timp3", "Agtr2", "MuSCs",3, 1, 36.46294533333333, 109.388836,         0.1676916294872044, 0.1676916294872044, 2, 0.6666666666666666, 1.093898333333333, 3.281695,  0.03053703709872749, 0.03053703709872749, 39.88675512855778,  358.98079615702,    0.005120805510796824, 0.005120805510796825),
    @("MuSCs","Timp3", "Agtr2", "ECs",  3, 1, 55.37138633333333, 166.114159,         0.2546507945619474, 0.2546507945619474, 1, 0.3333333333333333, 0.045282,           0.135846,  0.001264082841858775, 0.001264082841858776, 2.507327115946,     22.565944043514,    0.0003218997000714616,0.0003218997000714617),
    @("MuSCs","Timp3", "Agtr2", "FAPs", 3, 1, 55.37138633333333, 166.114159,         0.2546507945619474, 0.2546507945619474, 3, 1,                   34.682839,          104.048517, 0.9681988800594137,  0.9681988800594138,  1920.4368774058,    17283.9318966522,   0.2465526141011173,   0.2465526141011174),
    @("MuSCs","Timp3", "Agtr2", "MuSCs",3, 1, 55.37138633333333, 166.114159,         0.2546507945619474, 0.2546507945619474, 2, 0.6666666666666666, 1.093898333333333, 3.281695,  0.03053703709872749, 0.03053703709872749, 60.57066722438945,  545.136005019505,   0.007776280760758619, 0.007776280760758621)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}
